$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.614.01"
$ws.Range("E2").Value = "  -2.64%  "

# Row 3
$ws.Range("D3").Value = "2.324.73"
$ws.Range("E3").Value = "  -6.03%  "

# Row 4
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +1.24%  "

# Row 5
$ws.Range("D5").Value = "541.29"
$ws.Range("E5").Value = "  -2.01%  "

# Row 6
$ws.Range("D6").Value = "135.20"
$ws.Range("E6").Value = "  -8.52%  "

# Row 7
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8
$ws.Range("E8").Value = "  -11.41%  "

# Row 9
$ws.Range("D9").Value = "2.339.03"
$ws.Range("E9").Value = "  -5.51%  "

# Row 10
$ws.Range("E10").Value = "  -3.11%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("E12").Value = "  -3.02%  "

# Row 13
$ws.Range("E13").Value = "  -4.27%  "

# Row 14
$ws.Range("D14").Value = "24.36"
$ws.Range("E14").Value = "  -7.65%  "

# Row 15
$ws.Range("E15").Value = "  -5.42%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000159"
$ws.Range("E16").Value = "  -6.39%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "49.437.38"
$ws.Range("E17").Value = "  -20.46%  "

# Row 18
$ws.Range("D18").Value = "2.225.72"
$ws.Range("E18").Value = "  -9.77%  "

# Row 19
$ws.Range("D19").Value = "10.52"
$ws.Range("E19").Value = "  -4.60%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.06"
$ws.Range("E20").Value = "  -3.07%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "313.39"
$ws.Range("E21").Value = "  -2.76%  "

# Row 22
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  -8.17%  "

# Row 23
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.28%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "62.56"
$ws.Range("E24").Value = "  -2.42%  "

# Row 25
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").Value = "1.69"
$ws.Range("E25").Value = "  -10.42%  "

# Row 26
$ws.Range("D26").Value = "8.31"
$ws.Range("E26").Value = "  +6.21%  "

# Row 27
$ws.Range("E27").Value = "  -0.21%  "

# Row 28
$ws.Range("E28").Value = "  -5.50%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "7.84"
$ws.Range("E29").Value = "  -5.56%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "498.41"
$ws.Range("E30").Value = "  -7.90%  "

# Row 31
$ws.Range("D31").Value = "1.37"
$ws.Range("E31").Value = "  -8.63%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0864"
$ws.Range("E32").Value = "  -13.33%  "

# Row 33
$ws.Range("E33").Value = "  -2.83%  "

# Row 34
$ws.Range("E34").Value = "  -6.26%  "

# Row 35
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -7.79%  "

# Row 36
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.55%  "

# Row 37
$ws.Range("D37").Value = "4.57"
$ws.Range("E37").Value = "  -5.57%  "

# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "18.41"
$ws.Range("E38").Value = "  +0.27%  "

# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -3.09%  "

# Row 40
$ws.Range("D40").Value = "5.17"
$ws.Range("E40").Value = "  -10.96%  "

# Row 41
$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  -1.47%  "

# Row 42
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "140.56"
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").Value = "39.79"
$ws.Range("E44").Value = "  -1.87%  "

# Row 45
$ws.Range("D45").Value = "139.27"
$ws.Range("E45").Value = "  -4.26%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "3.52"
$ws.Range("E46").Value = "  -3.32%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  -11.39%  "

# Row 48
$ws.Range("D48").Value = "0.0508"
$ws.Range("E48").Value = "  -4.84%  "

# Row 49
$ws.Range("D49").Value = "19.06"
$ws.Range("E49").Value = "  -12.99%  "

# Row 50
$ws.Range("D50").Value = "0.565"
$ws.Range("E50").Value = "  -4.97%  "

# Row 51
$ws.Range("D51").Value = "0.0891"
$ws.Range("E51").Value = "  -4.86%  "
